$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Through 2022-08-19 -> Through 2022-08-20)
$ws.Name = "Through 2022-08-20"

# Update the header text in B1 (the running "through" date moved from Aug 19 to Aug 20)
$ws.Range("B1").Value = "August 2022 (through August 20)"

# New data column B holds the counts for 2022-08-20; existing rows bump up by
# the newly recorded carjacking(s) for that date, and a few new cells appear
# in other month columns for incidents recorded against 2022-08-28 data pull.

# Row 2
$ws.Range("B2").Value = 14
$ws.Range("AP2").Value = 3
$ws.Range("BF2").Value = 4

# Row 3
$ws.Range("R3").Value = 6

# Row 4
$ws.Range("B4").Value = 4
$ws.Range("J4").Value = 6

# Row 5
$ws.Range("B5").Value = 8

# Row 6
$ws.Range("BF6").Value = 3

# Row 9
$ws.Range("AP9").Value = 3

# Row 12
$ws.Range("B12").Value = 6

# Row 15
$ws.Range("B15").Value = 8
$ws.Range("J15").Value = 3

# Row 18
$ws.Range("B18").Value = 1

# Row 20
$ws.Range("AP20").Value = 3

# Row 26
$ws.Range("B26").Value = 3
$ws.Range("J26").Value = 2

# Row 30
$ws.Range("AP30").Value = 1
$ws.Range("AX30").Value = 1

# Row 31
$ws.Range("B31").Value = 2
$ws.Range("J31").Value = 1

# Row 34
$ws.Range("J34").Value = 1

# Row 37
$ws.Range("B37").Value = 5

# Row 42
$ws.Range("R42").Value = 1

# Row 46
$ws.Range("J46").Value = 4

# Row 50
$ws.Range("B50").Value = 3

# Row 53
$ws.Range("B53").Value = 2

# Row 64
$ws.Range("AP64").Value = 1

# Row 66
$ws.Range("B66").Value = 5
